# Generate Report for Archive
# - Replace "Ready for handoff" status text with "In Translation" everywhere
#   it appears (shared string used across Overview/zh-cn/de-de sheets).
# - Narrow the "Status" columns (E:F on Overview, C on zh-cn / de-de)
#   from ~17.22 to ~13.41 characters wide.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: put the string literal on the left of -eq; some cell values
        # come back as native booleans (e.g. the text "True") and
        # "$true -eq 'Ready for handoff'" would coerce the string to a
        # (non-empty => $true) boolean and false-positive match.
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
        }
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = 12.5
